$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6887f34fed985ac4c6bbde1587a351d31ee18e5/e2e/315556d7-10a4-49fb-9ee6-9762f0e7b5c0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f348b821050be3f8c538df4b52591468f68cf38e/e2e/315556d7-10a4-49fb-9ee6-9762f0e7b5c0.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P) to fit the new long message.
$wsZh.Range("P1").ColumnWidth = 39.2

# Latest Target File (I5): now handed off - becomes a hyperlink to the md file.
$i5 = $wsZh.Range("I5")
$i5.Value = "315556d7-10a4-49fb-9ee6-9762f0e7b5c0.md"
$wsZh.Hyperlinks.Add($i5, "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c6887f34fed985ac4c6bbde1587a351d31ee18e5/e2e/315556d7-10a4-49fb-9ee6-9762f0e7b5c0.md", "", "", "315556d7-10a4-49fb-9ee6-9762f0e7b5c0.md")

# Latest Handback File (J5)
$wsZh.Range("J5").Value = "315556d7-10a4-49fb-9ee6-9762f0e7b5c0.a990491c48a928a70a9e236dde36b6bfa04883b1.zh-cn.xlf"

# Latest Handback DateTime (K5)
$wsZh.Range("K5").Value = "2016-10-14 08:01:36"

# Error Detail (P5)
$wsZh.Range("P5").Value = $errorDetail

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

# Widen the "Error Detail" column (P) to fit the new long message.
$wsDe.Range("P1").ColumnWidth = 39.2

# Latest Target File (I5): now handed off - becomes a hyperlink to the md file.
$i5de = $wsDe.Range("I5")
$i5de.Value = "315556d7-10a4-49fb-9ee6-9762f0e7b5c0.md"
$wsDe.Hyperlinks.Add($i5de, "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c6887f34fed985ac4c6bbde1587a351d31ee18e5/e2e/315556d7-10a4-49fb-9ee6-9762f0e7b5c0.md", "", "", "315556d7-10a4-49fb-9ee6-9762f0e7b5c0.md")

# Latest Handback File (J5)
$wsDe.Range("J5").Value = "315556d7-10a4-49fb-9ee6-9762f0e7b5c0.a990491c48a928a70a9e236dde36b6bfa04883b1.de-de.xlf"

# Latest Handback DateTime (K5)
$wsDe.Range("K5").Value = "2016-10-14 08:01:53"

# Error Detail (P5)
$wsDe.Range("P5").Value = $errorDetail
